$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.775.13"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.483.22"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +4.70%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.27"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "658.29"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("E7").Value = "  +5.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.02"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.482.75"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  +4.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.39"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  +10.26%  "

$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.475.28"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.18"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +2.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.136.57"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +4.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000257"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.75"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.486.39"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +4.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.52"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +10.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.91"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +13.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.504"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -4.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "519.20"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +4.91%  "

$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000199"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.77"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  +4.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.57"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.53"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +4.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.658.18"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +4.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.26"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +13.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.82"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  +14.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("E33").Value = "  -2.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.188"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.590"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +7.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.94"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +9.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.86"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +4.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  +1.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.156"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +4.48%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "518.45"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.911"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +9.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.38"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("B45").Value = "ImmutableX"

$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +5.04%  "

$ws.Range("B46").Value = "VeChain"

$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0425"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +4.08%  "

$ws.Range("B47").Value = "Filecoin"

$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.66"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  +3.88%  "

$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("B49").Value = "dogwifhat"

$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.35"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  +7.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.22"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  +12.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.53"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  -0.91%  "
